# Update per-row Price (D) and Volume(1h) (E) values to refreshed figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.620.63"
$ws.Range("E2").Value = "  +2.13%  "

$ws.Range("D3").Value = "1.887.92"
$ws.Range("E3").Value = "  +0.24%  "

$ws.Range("E4").Value = "  +0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.18"
$ws.Range("E5").Value = "  +1.15%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4907"
$ws.Range("E7").Value = "  -0.68%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2942"
$ws.Range("E8").Value = "  -0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06762"
$ws.Range("E9").Value = "  +1.71%  "

$ws.Range("D10").Value = "1.888.22"
$ws.Range("E10").Value = "  +0.30%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.16"
$ws.Range("E11").Value = "  +2.30%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07228"
$ws.Range("E12").Value = "  +0.36%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "91.06"
$ws.Range("E13").Value = "  +5.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.057"
$ws.Range("E14").Value = "  +3.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6757"
$ws.Range("E15").Value = "  +0.96%  "

$ws.Range("D16").Value = "30.583.65"
$ws.Range("E16").Value = "  +2.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007956"
$ws.Range("E17").Value = "  +1.35%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("E18").Value = "  +0.20%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.17"
$ws.Range("E19").Value = "  +2.67%  "

$ws.Range("D20").Value = "2.132.79"
$ws.Range("E20").Value = "  +0.67%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("E21").Value = "  +0.11%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.812"
$ws.Range("E22").Value = "  +0.47%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "182.94"
$ws.Range("E23").Value = "  +28.14%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.060"
$ws.Range("E24").Value = "  +3.00%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.349"
$ws.Range("E25").Value = "  +2.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.38"
$ws.Range("E26").Value = "  +3.27%  "

$ws.Range("E27").Value = "  +11.19%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.901"
$ws.Range("E28").Value = "  -1.26%  "

$ws.Range("E29").Value = "  +0.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.314"
$ws.Range("E30").Value = "  +2.20%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09034"
$ws.Range("E31").Value = "  +2.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.997"
$ws.Range("E32").Value = "  -0.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05193"
$ws.Range("E33").Value = "  +2.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7539"
$ws.Range("E34").Value = "  +5.40%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.111"
$ws.Range("E35").Value = "  -0.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.750"
$ws.Range("E36").Value = "  +3.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01843"
$ws.Range("E37").Value = "  +2.58%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.660"
$ws.Range("E38").Value = "  -1.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.141"
$ws.Range("E39").Value = "  -1.52%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9360"
$ws.Range("E40").Value = "  +0.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4414"
$ws.Range("E41").Value = "  +4.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "105.46"
$ws.Range("E42").Value = "  +2.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.22%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.736"
$ws.Range("E44").Value = "  -0.53%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.590"
$ws.Range("E45").Value = "  +2.35%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1337"
$ws.Range("E46").Value = "  +5.01%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.696"
$ws.Range("E49").Value = "  +4.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3925"
$ws.Range("E50").Value = "  +3.95%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.51"
$ws.Range("E51").Value = "  +2.72%  "

# Rows 47/48: NEARProtocol and Cronos swap rank positions with refreshed values.
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.451"
$ws.Range("E47").Value = "  +8.22%  "

$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05844"
$ws.Range("E48").Value = "  +2.87%  "
